$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "57.324.77"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.098.24"
$ws.Range("E3").Value = "  -0.29%  "
Set-TextValue "D5" "524.16"
$ws.Range("E5").Value = "  +0.13%  "
Set-TextValue "D6" "136.96"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.097.49"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  +2.23%  "
Set-TextValue "D10" "7.22"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  -0.99%  "
Set-TextValue "D12" "0.395"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "3.637.41"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  +2.72%  "
Set-TextValue "D15" "25.21"
$ws.Range("E15").Value = "  -3.45%  "
Set-TextValue "D16" "0.0000163"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "57.424.48"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "3.104.37"
$ws.Range("E18").Value = "  -0.03%  "
Set-TextValue "D19" "5.93"
$ws.Range("E19").Value = "  -2.81%  "
Set-TextValue "D20" "12.46"
$ws.Range("E20").Value = "  -2.55%  "
Set-TextValue "D21" "7.87"
$ws.Range("E21").Value = "  -2.43%  "
Set-TextValue "D22" "346.73"
$ws.Range("E22").Value = "  +2.51%  "
Set-TextValue "D23" "5.78"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  +0.05%  "
Set-TextValue "D25" "68.19"
$ws.Range("E25").Value = "  +2.13%  "
Set-TextValue "D26" "0.502"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("E27").Value = "  -0.96%  "
Set-TextValue "D28" "0.997"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "0.0₃0905"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "7.36"
$ws.Range("E31").Value = "  +2.35%  "
Set-TextValue "D32" "1.87"
$ws.Range("E32").Value = "  +0.39%  "
Set-TextValue "D33" "6.00"
$ws.Range("E33").Value = "  -7.80%  "
Set-TextValue "D34" "20.96"
$ws.Range("E34").Value = "  -0.11%  "
Set-TextValue "D35" "1.16"
$ws.Range("E35").Value = "  -2.91%  "
Set-TextValue "D36" "4.91"
$ws.Range("E36").Value = "  +6.00%  "
Set-TextValue "D37" "158.37"
$ws.Range("E37").Value = "  +0.04%  "
Set-TextValue "D38" "6.13"
$ws.Range("E38").Value = "  +0.36%  "
Set-TextValue "D39" "25.89"
$ws.Range("E39").Value = "  -4.39%  "
Set-TextValue "D40" "1.23"
$ws.Range("E40").Value = "  -4.11%  "
Set-TextValue "D41" "4.18"
$ws.Range("E41").Value = "  +5.97%  "
Set-TextValue "D42" "0.0663"
$ws.Range("E42").Value = "  +0.61%  "
Set-TextValue "D43" "1.60"
$ws.Range("E43").Value = "  +5.68%  "
Set-TextValue "D44" "0.697"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").Value = "3.142.66"
$ws.Range("E45").Value = "  -0.25%  "
Set-TextValue "D46" "36.40"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.340.50"
$ws.Range("E48").Value = "  +1.52%  "
Set-TextValue "D49" "0.0267"
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "6.02"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D51" "0.949"
$ws.Range("E51").Value = "  -2.84%  "
